# Insert a new weekly price-observation row at row 436, pushing the
# existing rows 436:534 down to 437:535 (dimension grows from R534 to R535).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(436).Insert()

# Columns A,B,C,E,F,G,H,I,R are constant for every data row in this
# "Mercado" sheet (Feria Lagunitas de Puerto Montt / Cilantro), so copy
# them straight from the row directly below (the former row 436, now 437).
$ws.Cells.Item(436, 1).Value2 = $ws.Cells.Item(437, 1).Value2   # Mercado ID
$ws.Cells.Item(436, 2).Value2 = $ws.Cells.Item(437, 2).Value2   # Mercado
$ws.Cells.Item(436, 3).Value2 = $ws.Cells.Item(437, 3).Value2   # Region
$ws.Cells.Item(436, 5).Value2 = $ws.Cells.Item(437, 5).Value2   # Codreg
$ws.Cells.Item(436, 6).Value2 = $ws.Cells.Item(437, 6).Value2   # Categoria ID
$ws.Cells.Item(436, 7).Value2 = $ws.Cells.Item(437, 7).Value2   # Categoria
$ws.Cells.Item(436, 8).Value2 = $ws.Cells.Item(437, 8).Value2   # Variedad
$ws.Cells.Item(436, 9).Value2 = $ws.Cells.Item(437, 9).Value2   # Calidad
$ws.Cells.Item(436, 18).Value2 = $ws.Cells.Item(437, 18).Value2 # Clasificacion

# New observation's own data.
$ws.Cells.Item(436, 4).Value2 = 45204                       # Fecha
$ws.Cells.Item(436, 10).Value2 = 120                        # Volumen
$ws.Cells.Item(436, 11).Value2 = 13000                      # Precio minimo
$ws.Cells.Item(436, 12).Value2 = 13000                      # Precio maximo
$ws.Cells.Item(436, 13).Value2 = 13000                      # Precio promedio ponderado
$ws.Cells.Item(436, 14).Value2 = "$/caja 36 atados"         # Unidad de comercializacion
$ws.Cells.Item(436, 15).Value2 = "Región Metropolitana"     # Origen
$ws.Cells.Item(436, 16).Value2 = 361                        # Precio $/Kg
$ws.Cells.Item(436, 17).Value2 = 36                         # Kg o Unidades
